# Auto-generated edit script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "80.907.46"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "3.134.67"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.39"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.20"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.279"
$ws.Range("E7").Value = "  +23.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "3.129.80"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.570"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000249"
$ws.Range("E12").Value = "  +11.34%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "3.709.80"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "80.830.48"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").Value = "3.110.87"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("E19").Value = "  +9.01%  "
$ws.Range("E20").Value = "  -4.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "427.84"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.88"
$ws.Range("E22").Value = "  -4.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.04"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.11"
$ws.Range("E25").Value = "  +8.19%  "
$ws.Range("D26").Value = "3.287.53"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "75.37"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.71"
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.85"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "552.50"
$ws.Range("E33").Value = "  +7.72%  "
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.151"
$ws.Range("E35").Value = "  +12.79%  "
$ws.Range("B36").Value = "Cronos"
$ws.Range("C36").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.137"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.40"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.96"
$ws.Range("E40").Value = "  +10.51%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.402"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("E42").Value = "  +3.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("E43").Value = "  +12.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.99"
$ws.Range("E44").Value = "  +20.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "158.82"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "186.00"
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.50"
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.758"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.46"
$ws.Range("E51").Value = "  +3.03%  "

Write-Output "Applied cryptos list update: $($wb.Name)"
